$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column A that currently read "All" should be changed to "Combined"
$rows = @(2, 5, 8, 11, 14, 17)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "All") {
        $cell.Value = "Combined"
    }
}
